# Commit: Common: Moved DTO resolution to another place
#
# The DTO resolution values "Another Value 001", "Va-Va-Value 001" and
# "Bad one" are no longer produced, so the cells that held them are
# cleared out entirely (not just blanked) on the "Nějaký import" sheet:
#   - B2/C2 ("Another Value 001" / "Va-Va-Value 001")
#   - A22, the trailing row  ("Bad one")
# This shrinks the sheet's used range from A1:C22 down to A1:C21 and -
# since nothing else references those three strings - the workbook's
# shared-string table loses exactly those three entries on save, which
# is also why every other sheet's <v> shared-string indices shift down.
#
# The view/selection state also moves: "Nějaký import" becomes the
# active/selected sheet (it used to be "tabs"), with a fresh selection.

$wb = $excel.ActiveWorkbook

$wsTabs   = $wb.Worksheets.Item("tabs")
$wsImport = $wb.Worksheets.Item("Nějaký import")

# Drop the cells/row that referenced the removed DTO values.
[void]$wsImport.Range("B2").ClearContents()
[void]$wsImport.Range("C2").ClearContents()
[void]$wsImport.Range("A22").ClearContents()

# "tabs" keeps its own remembered selection, just moved to B3, and is no
# longer the selected tab once another sheet is activated below.
[void]$wsTabs.Activate()
[void]$wsTabs.Range("B3").Select()

# "Nějaký import" becomes the active sheet with selection on D6.
[void]$wsImport.Activate()
[void]$wsImport.Range("D6").Select()
